$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '29.302.60'
$ws.Range('E2').Value2 = '  +0.59%  '
$ws.Range('D3').Value2 = '1.933.67'
$ws.Range('E3').Value2 = '  +1.50%  '
$ws.Range('D4').Value2 = '''1.003'
$ws.Range('E4').Value2 = '  +0.41%  '
$ws.Range('D5').Value2 = '''325.66'
$ws.Range('E5').Value2 = '  -0.11%  '
$ws.Range('D6').Value2 = '''1.001'
$ws.Range('E6').Value2 = '  +0.10%  '
$ws.Range('D7').Value2 = '''0.4622'
$ws.Range('E7').Value2 = '  +0.21%  '
$ws.Range('D8').Value2 = '''0.3875'
$ws.Range('E8').Value2 = '  -0.36%  '
$ws.Range('B9').Value2 = 'OKB'
$ws.Range('C9').Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value2 = '''45.90'
$ws.Range('E9').Value2 = '  -0.11%  '
$ws.Range('B10').Value2 = 'Dogecoin'
$ws.Range('C10').Value2 = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value2 = '''0.07815'
$ws.Range('E10').Value2 = '  -0.70%  '
$ws.Range('B11').Value2 = 'Polygon'
$ws.Range('C11').Value2 = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').Value2 = '''0.9751'
$ws.Range('E11').Value2 = '  -1.57%  '
$ws.Range('B12').Value2 = 'Solana'
$ws.Range('C12').Value2 = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').Value2 = '''22.65'
$ws.Range('E12').Value2 = '  +3.08%  '
$ws.Range('B13').Value2 = 'WrappedEther'
$ws.Range('C13').Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value2 = '1.918.18'
$ws.Range('E13').Value2 = '  +0.86%  '
$ws.Range('B14').Value2 = 'Polkadot'
$ws.Range('C14').Value2 = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value2 = '''5.795'
$ws.Range('E14').Value2 = '  +0.71%  '
$ws.Range('B15').Value2 = 'Chainlink'
$ws.Range('C15').Value2 = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value2 = '''7.083'
$ws.Range('E15').Value2 = '  +0.63%  '
$ws.Range('B16').Value2 = 'TRON'
$ws.Range('C16').Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value2 = '''0.07056'
$ws.Range('E16').Value2 = '  +0.54%  '
$ws.Range('B17').Value2 = 'Litecoin'
$ws.Range('C17').Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value2 = '''86.81'
$ws.Range('E17').Value2 = '  -1.48%  '
$ws.Range('B18').Value2 = 'BinanceUSD'
$ws.Range('C18').Value2 = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').Value2 = '''1.003'
$ws.Range('E18').Value2 = '  -0.03%  '
$ws.Range('B19').Value2 = 'ShibaInu'
$ws.Range('C19').Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value2 = '''0.000009743'
$ws.Range('E19').Value2 = '  -2.06%  '
$ws.Range('B20').Value2 = 'Avalanche'
$ws.Range('C20').Value2 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value2 = '''17.01'
$ws.Range('E20').Value2 = '  -0.49%  '
$ws.Range('B21').Value2 = 'Dai'
$ws.Range('C21').Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value2 = '''1.001'
$ws.Range('E21').Value2 = '  +0.05%  '
$ws.Range('B22').Value2 = 'WrappedBTC'
$ws.Range('C22').Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').Value2 = '29.291.02'
$ws.Range('E22').Value2 = '  +0.45%  '
$ws.Range('B23').Value2 = 'Uniswap'
$ws.Range('C23').Value2 = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value2 = '''5.475'
$ws.Range('E23').Value2 = '  +2.91%  '
$ws.Range('B24').Value2 = 'Cosmos'
$ws.Range('C24').Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value2 = '''11.06'
$ws.Range('E24').Value2 = '  -0.54%  '
$ws.Range('B25').Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value2 = '2.163.06'
$ws.Range('E25').Value2 = '  +3.89%  '
$ws.Range('B26').Value2 = 'Toncoin'
$ws.Range('C26').Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value2 = '''2.094'
$ws.Range('E26').Value2 = '  +0.18%  '
$ws.Range('B27').Value2 = 'Monero'
$ws.Range('C27').Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value2 = '''157.51'
$ws.Range('E27').Value2 = '  +0.65%  '
$ws.Range('B28').Value2 = 'EthereumClassic'
$ws.Range('C28').Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value2 = '''19.37'
$ws.Range('E28').Value2 = '  -0.51%  '
$ws.Range('B29').Value2 = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value2 = '''5.762'
$ws.Range('E29').Value2 = '  -2.22%  '
$ws.Range('B30').Value2 = 'BitcoinCash'
$ws.Range('C30').Value2 = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value2 = '''118.85'
$ws.Range('E30').Value2 = '  +0.07%  '
$ws.Range('B31').Value2 = 'LidoDAOToken'
$ws.Range('C31').Value2 = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D31').Value2 = '''1.845'
$ws.Range('E31').Value2 = '  -1.74%  '
$ws.Range('B32').Value2 = 'Stellar'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value2 = '''0.09331'
$ws.Range('E32').Value2 = '  -0.18%  '
$ws.Range('B33').Value2 = 'ImmutableX'
$ws.Range('C33').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value2 = '''0.8639'
$ws.Range('E33').Value2 = '  -3.57%  '
$ws.Range('B34').Value2 = 'Filecoin'
$ws.Range('C34').Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value2 = '''5.167'
$ws.Range('E34').Value2 = '  -1.21%  '
$ws.Range('B35').Value2 = 'ARBITRUM'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value2 = '''1.301'
$ws.Range('E35').Value2 = '  -1.55%  '
$ws.Range('B36').Value2 = 'HuobiToken'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value2 = '''3.080'
$ws.Range('E36').Value2 = '  -1.91%  '
$ws.Range('B37').Value2 = 'Hedera'
$ws.Range('C37').Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value2 = '''0.05776'
$ws.Range('E37').Value2 = '  -0.14%  '
$ws.Range('B38').Value2 = 'TrustWalletToken'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value2 = '''1.157'
$ws.Range('E38').Value2 = '  -1.14%  '
$ws.Range('B39').Value2 = 'VeChain'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value2 = '''0.02081'
$ws.Range('E39').Value2 = '  -0.25%  '
$ws.Range('B40').Value2 = 'FraxShare'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value2 = '''7.635'
$ws.Range('E40').Value2 = '  -0.51%  '
$ws.Range('B41').Value2 = 'TheSandbox'
$ws.Range('C41').Value2 = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value2 = '''0.5657'
$ws.Range('E41').Value2 = '  -0.65%  '
$ws.Range('D42').Value2 = '''0.1779'
$ws.Range('E42').Value2 = '  -1.63%  '
$ws.Range('B43').Value2 = 'PEPE'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value2 = '''0.000003051'
$ws.Range('E43').Value2 = '  +22.71%  '
$ws.Range('B44').Value2 = 'Aptos'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value2 = '''9.392'
$ws.Range('E44').Value2 = '  -3.55%  '
$ws.Range('B45').Value2 = 'MXToken'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value2 = '''2.713'
$ws.Range('E45').Value2 = '  +6.43%  '
$ws.Range('B46').Value2 = 'Decentraland'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value2 = '''0.5263'
$ws.Range('E46').Value2 = '  -1.66%  '
$ws.Range('B47').Value2 = 'EnergySwap'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value2 = '''11.40'
$ws.Range('E47').Value2 = '  -3.53%  '
$ws.Range('B48').Value2 = 'Cronos'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value2 = '''0.06866'
$ws.Range('E48').Value2 = '  -1.82%  '
$ws.Range('B49').Value2 = 'RenderToken'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value2 = '''2.073'
$ws.Range('E49').Value2 = '  -4.40%  '
$ws.Range('B50').Value2 = 'NEARProtocol'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value2 = '''1.814'
$ws.Range('E50').Value2 = '  -1.43%  '
$ws.Range('B51').Value2 = 'Quant'
$ws.Range('C51').Value2 = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value2 = '''111.16'
$ws.Range('E51').Value2 = '  -1.91%  '
